# Unificação de alguns testes e pequenas refatorações
#
# Planilha2 (A2:B3) already holds the "Teste_BuscaLupa" / "Teste_BuscaLupaFalha"
# test rows. Bring copies of them into Planilha1 so the suite lives in one
# place, widen the column that now holds the longer product names, rename a
# stray generated username, and leave the selections where the user ended up.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# --- Planilha1: copy in the two "BuscaLupa" rows from Planilha2 -------------
$ws1.Range("A5").Value2 = $ws2.Range("A2").Value2
$ws1.Range("B5").Value2 = $ws2.Range("B2").Value2

$ws1.Range("A6").Value2 = $ws2.Range("A3").Value2
$ws1.Range("B6").Value2 = $ws2.Range("B3").Value2

# Match the look of the other data rows (thin box border, same as A2:B3).
$ws1.Range("A2:B3").Copy()
$ws1.Range("A5:B6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Keep row 9's empty helper cell styled like C6 (border-less, underlined font).
$ws1.Range("C6").Copy()
$ws1.Range("B9").PasteSpecial(-4122)      # xlPasteFormats
$excel.CutCopyMode = 0

# Rename the stray generated test id referenced from B2.
$ws1.Range("B2").Value2 = "usertests3884"

# Column B needs to be wide enough for the newly-added product names.
$ws1.Columns.Item(2).ColumnWidth = 24.14

# Selections, reflecting where the user ended up after the edit.
$ws2.Activate()
$ws2.Range("A2:B3").Select()

$ws1.Activate()
$ws1.Range("C8").Select()
